# Form B - 2/2/23 Minutes and Plan Form
# Fill in the post-meeting minutes & plan placeholders, and drop the
# unused placeholder bullets / spacer paragraph that are no longer needed.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace the placeholder "." bullets with the real minutes / plan text.
#    (Paragraph indices are stable here since we only touch Range.Text,
#    which does not add or remove paragraphs.)
# ---------------------------------------------------------------------------

$d.Paragraphs.Item(23).Range.Text = "Regarding the USB library issue: Use the USB contributors full repository. Continue to try working on this component, try communicating with the contributor to find help in solving any problems."
$d.Paragraphs.Item(24).Range.Text = "Discussion on dissertation writing, highlighting new sections and their placement within the report. "
$d.Paragraphs.Item(25).Range.Text = "Discussion of new content that could be added to report to enrich existing sections."
$d.Paragraphs.Item(26).Range.Text = "Brief discussion on expectations of presentation."

$d.Paragraphs.Item(39).Range.Text = "Presentation writing."
$d.Paragraphs.Item(40).Range.Text = "Continuing work and investigation into the Rust for Linux USB contributions. "

$d.Paragraphs.Item(46).Range.Text = "Continue work on writing for dissertation/final report."

# ---------------------------------------------------------------------------
# 2) Drop the now-unused placeholder bullets (and one spare blank spacer
#    paragraph). Deletions are performed from the bottom of the document
#    upward so that earlier paragraph indices remain valid.
# ---------------------------------------------------------------------------

# "Beyond the next month" list (numId 3): keep only the first item.
$d.Paragraphs.Item(49).Range.Delete()
$d.Paragraphs.Item(48).Range.Delete()
$d.Paragraphs.Item(47).Range.Delete()

# "For the next month" list (numId 2): keep only the first two items.
$d.Paragraphs.Item(42).Range.Delete()
$d.Paragraphs.Item(41).Range.Delete()

# Spacer paragraph right after the Minutes list.
$d.Paragraphs.Item(29).Range.Delete()

# Minutes list (numId 1): keep only the first four items.
$d.Paragraphs.Item(28).Range.Delete()
$d.Paragraphs.Item(27).Range.Delete()

Write-Output "Paragraph count after edits: $($d.Paragraphs.Count)"
